$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
for ($i=1; $i -le 12; $i++) {
    $c = $tcs.Item($i)
    Write-Output "$i Type=$($c.Type)"
}
